$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Font.Bold = $true
Write-Host $ws.Range("F2").Font.Bold
